# Update the macro correlation table (WorkingFolder/Tables/macro_corr.xlsx)
# with refreshed values exported from Stata.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the old values first so that any shared strings that become
# unused are dropped when the workbook is re-saved.
$ws.Range("B2:N7").ClearContents() | Out-Null

# New correlation values (row 2 = median:var ... row 7 = mean:rvar,
# columns B..N). Values are written as text (matching the source data,
# which includes significance markers such as "*" and "**").
$newValues = [ordered]@{
    "B2" = "-0.11"
    "C2" = "-0.12"
    "D2" = "-0.2*"
    "E2" = "-0.14"
    "F2" = "-0.05"
    "G2" = "-0.06"
    "H2" = "-0.09"
    "I2" = "-0.17"
    "J2" = "-0.08"
    "K2" = "0.0"
    "L2" = "0.09"
    "M2" = "0.01"
    "N2" = "0.12"
    "B3" = "-0.11"
    "C3" = "-0.11"
    "D3" = "-0.21*"
    "E3" = "-0.12"
    "F3" = "0.01"
    "G3" = "-0.03"
    "H3" = "-0.09"
    "I3" = "-0.21*"
    "J3" = "-0.11"
    "K3" = "0.01"
    "L3" = "0.1"
    "M3" = "0.03"
    "N3" = "0.11"
    "B4" = "-0.08"
    "C4" = "0.04"
    "D4" = "0.06"
    "E4" = "0.01"
    "F4" = "-0.15"
    "G4" = "-0.17"
    "H4" = "-0.21*"
    "I4" = "-0.28**"
    "J4" = "-0.22*"
    "K4" = "-0.05"
    "L4" = "0.01"
    "M4" = "0.02"
    "N4" = "-0.05"
    "B5" = "-0.08"
    "C5" = "0.01"
    "D5" = "0.01"
    "E5" = "-0.07"
    "F5" = "-0.21*"
    "G5" = "-0.18"
    "H5" = "-0.12"
    "I5" = "-0.09"
    "J5" = "-0.14"
    "K5" = "-0.22*"
    "L5" = "-0.25**"
    "M5" = "-0.06"
    "N5" = "-0.04"
    "B6" = "-0.07"
    "C6" = "0.03"
    "D6" = "-0.01"
    "E6" = "-0.09"
    "F6" = "-0.22*"
    "G6" = "-0.18"
    "H6" = "-0.13"
    "I6" = "-0.15"
    "J6" = "-0.2*"
    "K6" = "-0.23**"
    "L6" = "-0.22*"
    "M6" = "-0.09"
    "N6" = "-0.05"
    "B7" = "-0.08"
    "C7" = "0.07"
    "D7" = "0.19*"
    "E7" = "0.15"
    "F7" = "-0.16"
    "G7" = "-0.23**"
    "H7" = "-0.17"
    "I7" = "-0.08"
    "J7" = "-0.12"
    "K7" = "-0.17"
    "L7" = "-0.07"
    "M7" = "0.01"
    "N7" = "-0.04"
}

foreach ($addr in $newValues.Keys) {
    $text = $newValues[$addr]
    $cell = $ws.Range($addr)

    # Writing the value directly would make Excel store purely-numeric-looking
    # strings (e.g. '-0.11') as numbers. The source workbook stores every value
    # in this table as text, so build it via a text formula and then convert the
    # formula to a static value (Copy + PasteSpecial values) to force a text cell
    # without touching the cell style.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

$excel.CutCopyMode = 0

Write-Host "macro_corr.xlsx correlation table updated"
